$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 46070
$ws.Range("C3").Value = 46070
$ws.Range("C4").Value = 46070
$ws.Range("C5").Value = 46070
$ws.Range("C6").Value = 46070
$ws.Range("C7").Value = 46070
$ws.Range("C8").Value = 46070
$ws.Range("C9").Value = 46070
$ws.Range("C10").Value = 46070
$ws.Range("C11").Value = 46070
$ws.Range("C12").Value = 46070
$ws.Range("C13").Value = 46070
$ws.Range("C14").Value = 46070
$ws.Range("C15").Value = 46070
$ws.Range("C16").Value = 46070
$ws.Range("C17").Value = 46070
$ws.Range("C18").Value = 46070
$ws.Range("C19").Value = 46070
$ws.Range("C20").Value = 46070
$ws.Range("C21").Value = 46070
$ws.Range("C22").Value = 46070
$ws.Range("C23").Value = 46070
$ws.Range("C24").Value = 46070
$ws.Range("A25").Value = "A 63378-2025"
$ws.Range("B25").Value = 46010.72528935185
$ws.Range("C25").Value = 46070
$ws.Range("G25").Value = 4.8
$ws.Range("A26").Value = "A 34073-2025"
$ws.Range("B26").Value = 45845
$ws.Range("C26").Value = 46070
$ws.Range("G26").Value = 2.9
$ws.Range("A27").Value = "A 35047-2025"
$ws.Range("B27").Value = 45852.35094907408
$ws.Range("C27").Value = 46070
$ws.Range("G27").Value = 4.2
$ws.Range("A28").Value = "A 37245-2025"
$ws.Range("B28").Value = 45876.39396990741
$ws.Range("C28").Value = 46070
$ws.Range("G28").Value = 1.1
$ws.Range("A29").Value = "A 2533-2026"
$ws.Range("B29").Value = 46037.44622685185
$ws.Range("C29").Value = 46070
$ws.Range("G29").Value = 1.4
$ws.Range("A30").Value = "A 2535-2026"
$ws.Range("B30").Value = 46037.4490625
$ws.Range("C30").Value = 46070
$ws.Range("G30").Value = 1.6
$ws.Range("A31").Value = "A 55724-2024"
$ws.Range("B31").Value = 45622.81020833334
$ws.Range("C31").Value = 46070
$ws.Range("G31").Value = 0.5
$ws.Range("A32").Value = "A 35955-2023"
$ws.Range("B32").Value = 45148
$ws.Range("C32").Value = 46070
$ws.Range("G32").Value = 1.1
$ws.Range("A33").Value = "A 61217-2024"
$ws.Range("B33").Value = 45645.65825231482
$ws.Range("C33").Value = 46070
$ws.Range("G33").Value = 0.8
$ws.Range("A34").Value = "A 56002-2024"
$ws.Range("B34").Value = 45623.7328587963
$ws.Range("C34").Value = 46070
$ws.Range("G34").Value = 3.4
$ws.Range("A35").Value = "A 60728-2022"
$ws.Range("B35").Value = 44913
$ws.Range("C35").Value = 46070
$ws.Range("G35").Value = 1.2
$ws.Range("A36").Value = "A 19935-2025"
$ws.Range("B36").Value = 45771
$ws.Range("C36").Value = 46070
$ws.Range("G36").Value = 2.1
$ws.Range("A37").Value = "A 3222-2022"
$ws.Range("B37").Value = 44582
$ws.Range("C37").Value = 46070
$ws.Range("G37").Value = 1.8
$ws.Range("C38").Value = 46070
$ws.Range("C39").Value = 46070
